$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = "2025/12/04 06:00"
$ws.Range("B38").Value = "32,744位本"
$ws.Range("C38").Value = "87位 広告・宣伝 (本)"
$ws.Range("D38").Value = "140位商業デザイン"
$ws.Range("E38").Value = "1,749位ビジネス実用本"
$ws.Range("F38").Value = "-"
$ws.Range("G38").Value = "-"
